$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 311; everything from old row 311 onward
# shifts down by one (old 311 -> 312, ..., old 379 -> 380).
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new data record.
$ws.Cells.Item(311, 1).Value = 8
$ws.Cells.Item(311, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(311, 3).Value = "Coquimbo"
$ws.Cells.Item(311, 4).Value = 44641
$ws.Cells.Item(311, 5).Value = 4
$ws.Cells.Item(311, 6).Value = 100114001
$ws.Cells.Item(311, 7).Value = "Papa"
$ws.Cells.Item(311, 8).Value = "Asterix"
$ws.Cells.Item(311, 9).Value = "1a (cosecha)"
$ws.Cells.Item(311, 10).Value = 2400
$ws.Cells.Item(311, 11).Value = 8000
$ws.Cells.Item(311, 12).Value = 9000
$ws.Cells.Item(311, 13).Value = 8500
$ws.Cells.Item(311, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(311, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(311, 16).Value = 340
$ws.Cells.Item(311, 17).Value = 25
$ws.Cells.Item(311, 18).Value = "Hortaliza"
